$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 7053.3125
$ws.Range("J51").Value = 7543.077
$ws.Range("L51").Value = 7543.077
$ws.Range("N51").Value = -8511.077000000001
$ws.Range("H116").Value = 6793.375
$ws.Range("I116").Value = 6891.1665
$ws.Range("J116").Value = 6500
$ws.Range("K116").Value = 6891.1665
$ws.Range("L116").Value = 6500
$ws.Range("M116").Value = -3449.1665
$ws.Range("N116").Value = -13384
$ws.Range("H132").Value = 3939.4092
$ws.Range("I132").Value = 3456.4211
$ws.Range("K132").Value = 10369.2633
$ws.Range("M132").Value = -7839.263300000001
$ws.Range("H137").Value = 2198.2
$ws.Range("I137").Value = 892.6923
$ws.Range("K137").Value = 2678.0769
$ws.Range("M137").Value = -128.0769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 17500
$ws.Range("J24").Value = 17500
$ws.Range("L24").Value = 17500
$ws.Range("N24").Value = -18248
$ws.Range("H45").Value = 1721.4286
$ws.Range("I45").Value = 1607.5714
$ws.Range("K45").Value = 1607.5714
$ws.Range("M45").Value = -1230.5714
$ws.Range("H92").Value = 75785
$ws.Range("J92").Value = 75785
$ws.Range("L92").Value = 75785
$ws.Range("N92").Value = -80777
$ws.Range("H96").Value = 18500
$ws.Range("J96").Value = 18500
$ws.Range("L96").Value = 18500
$ws.Range("N96").Value = -23992
$ws.Range("H100").Value = 17500
$ws.Range("J100").Value = 17500
$ws.Range("L100").Value = 17500
$ws.Range("N100").Value = -19664
$ws.Range("H132").Value = 3034.9524
$ws.Range("I132").Value = 2714.1177
$ws.Range("K132").Value = 8142.353099999999
$ws.Range("M132").Value = -5612.353099999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2194.5881
$ws.Range("I105").Value = 1068.1666
$ws.Range("J105").Value = 4898
$ws.Range("K105").Value = 1068.1666
$ws.Range("L105").Value = 4898
$ws.Range("M105").Value = 678.8334
$ws.Range("N105").Value = -8392
$ws.Range("H137").Value = 63250
$ws.Range("J137").Value = 63250
$ws.Range("L137").Value = 63250
$ws.Range("N137").Value = -73450

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 6000
$ws.Range("J45").Value = 10000
$ws.Range("L45").Value = 30000
$ws.Range("N45").Value = -31064
$ws.Range("H68").Value = 4122.857
$ws.Range("I68").Value = 4290.6665
$ws.Range("K68").Value = 12871.9995
$ws.Range("M68").Value = -12060.9995
$ws.Range("H71").Value = 4122.857
$ws.Range("I71").Value = 4290.6665
$ws.Range("K71").Value = 38615.9985
$ws.Range("M71").Value = -34559.9985
$ws.Range("H107").Value = 996.2105
$ws.Range("J107").Value = 1112.1428
$ws.Range("L107").Value = 3336.4284
$ws.Range("N107").Value = -7176.428400000001
$ws.Range("H131").Value = 24677.812
$ws.Range("I131").Value = 140002.25
$ws.Range("K131").Value = 420006.75
$ws.Range("M131").Value = -414966.75
$ws.Range("H132").Value = 5190.625
$ws.Range("J132").Value = 5790.7144
$ws.Range("L132").Value = 52116.4296
$ws.Range("N132").Value = -57176.4296
$ws.Range("H138").Value = 1653.7142
$ws.Range("I138").Value = 1653.7142
$ws.Range("K138").Value = 4961.142599999999
$ws.Range("M138").Value = 178.8574000000008
$ws.Range("H139").Value = 4423.727
$ws.Range("I139").Value = 2701.2104
$ws.Range("J139").Value = 15333
$ws.Range("K139").Value = 8103.6312
$ws.Range("L139").Value = 45999
$ws.Range("M139").Value = -2963.6312
$ws.Range("N139").Value = -56279

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 22400
$ws.Range("J18").Value = 30000
$ws.Range("L18").Value = 30000
$ws.Range("N18").Value = -30586
$ws.Range("H43").Value = 1610.6666
$ws.Range("I43").Value = 1610.6666
$ws.Range("K43").Value = 1610.6666
$ws.Range("M43").Value = -1459.6666
$ws.Range("H70").Value = 7965.6924
$ws.Range("I70").Value = 5193.8
$ws.Range("J70").Value = 9698.125
$ws.Range("K70").Value = 5193.8
$ws.Range("L70").Value = 9698.125
$ws.Range("M70").Value = -4923.8
$ws.Range("N70").Value = -10238.125
$ws.Range("H73").Value = 7965.6924
$ws.Range("I73").Value = 5193.8
$ws.Range("J73").Value = 9698.125
$ws.Range("K73").Value = 5193.8
$ws.Range("L73").Value = 9698.125
$ws.Range("M73").Value = -4257.8
$ws.Range("N73").Value = -11570.125
$ws.Range("H80").Value = 6709.091
$ws.Range("I80").Value = 4172
$ws.Range("J80").Value = 9753.6
$ws.Range("K80").Value = 4172
$ws.Range("L80").Value = 9753.6
$ws.Range("M80").Value = -3174
$ws.Range("N80").Value = -11749.6
$ws.Range("H83").Value = 6709.091
$ws.Range("I83").Value = 4172
$ws.Range("J83").Value = 9753.6
$ws.Range("K83").Value = 20860
$ws.Range("L83").Value = 48768
$ws.Range("M83").Value = -15868
$ws.Range("N83").Value = -58752
$ws.Range("H113").Value = 2363.5
$ws.Range("I113").Value = 1899.3334
$ws.Range("K113").Value = 1899.3334
$ws.Range("M113").Value = 270.6666
$ws.Range("H132").Value = 3384.65
$ws.Range("J132").Value = 4275.8
$ws.Range("L132").Value = 12827.4
$ws.Range("N132").Value = -17887.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 417.29413
$ws.Range("I55").Value = 349.625
$ws.Range("K55").Value = 349.625
$ws.Range("M55").Value = -176.625
$ws.Range("H68").Value = 2473.7334
$ws.Range("I68").Value = 2446.0908
$ws.Range("K68").Value = 2446.0908
$ws.Range("M68").Value = -1697.0908
$ws.Range("H71").Value = 2473.7334
$ws.Range("I71").Value = 2446.0908
$ws.Range("K71").Value = 12230.454
$ws.Range("M71").Value = -8486.454
$ws.Range("H104").Value = 28837.5
$ws.Range("J104").Value = 28837.5
$ws.Range("L104").Value = 28837.5
$ws.Range("N104").Value = -35825.5
$ws.Range("H132").Value = 3867.875
$ws.Range("I132").Value = 3695.6667
$ws.Range("J132").Value = 4040.0833
$ws.Range("K132").Value = 11087.0001
$ws.Range("L132").Value = 12120.2499
$ws.Range("M132").Value = -8557.000100000001
$ws.Range("N132").Value = -17180.2499

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 14370
$ws.Range("J104").Value = 14370
$ws.Range("L104").Value = 14370
$ws.Range("N104").Value = -21358
$ws.Range("H125").Value = 17895
$ws.Range("J125").Value = 17895
$ws.Range("L125").Value = 17895
$ws.Range("N125").Value = -27735
$ws.Range("H132").Value = 3178.138
$ws.Range("I132").Value = 3113.0715
$ws.Range("K132").Value = 9339.2145
$ws.Range("M132").Value = -6809.2145
$ws.Range("H136").Value = 3505.2144
$ws.Range("I136").Value = 2197.923
$ws.Range("K136").Value = 6593.768999999999
$ws.Range("M136").Value = -4043.768999999999
